$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I header
$ws.Range("I1").Value = "Other found locations"

# Row 2 - update Authors (E2) and add Other found locations (I2)
$ws.Range("E2").Value = "[Martin%Chalumeau%NULL%1,  Philippe%Bidet%NULL%1,  Gérard%Lina%NULL%1,  Mostapha%Mokhtari%NULL%1,  Marie-Claude%André%NULL%1,  Dominique%Gendrel%NULL%1,  Edouard%Bingen%NULL%1,  Josette%Raymond%j.raymond@svp.ap-hop-paris.fr%1]"
$ws.Range("I2").Value = "_PMC"

# Row 3
$ws.Range("E3").Value = "[Hae-Sung%Nam%NULL%0,  Mi-Yeon%Yeon%NULL%2,  Mi-Yeon%Yeon%NULL%0,  Jung Wan%Park%NULL%2,  Jung Wan%Park%NULL%0,  Jee-Young%Hong%NULL%2,  Jee-Young%Hong%NULL%0,  Ji Woong%Son%NULL%2,  Ji Woong%Son%NULL%0]"
$ws.Range("I3").Value = "_PMC"

# Row 4
$ws.Range("E4").Value = "[Michael D.%Christian%NULL%1,  Mona%Loutfy%NULL%1,  L. Clifford%McDonald%NULL%2,  Kenneth F.%Martinez%NULL%1,  Mariana%Ofner%NULL%1,  Tom%Wong%NULL%1,  Tamara%Wallington%NULL%1,  Wayne L.%Gold%NULL%1,  Barbara%Mederski%NULL%1,  Karen%Green%NULL%3,  Donald E.%Low%NULL%1,  NULL%NULL%NULL%15]"
$ws.Range("I4").Value = "_PMC"

# Row 5 - only adds I5 (empty string)
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("E6").Value = "[ J.%Knapp%null%2,   M.A.%Weigand%null%1,   E.%Popp%null%1,  J.%Knapp%null%0,  M.A.%Weigand%null%1,  E.%Popp%null%1]"
$ws.Range("I6").Value = ""

# Row 7
$ws.Range("E7").Value = "[Wei%Liu%NULL%0,  Fang%Tang%NULL%1,  Li‐Qun%Fang%NULL%1,  Sake J.%De Vlas%NULL%1,  Huai‐Jian%Ma%NULL%1,  Jie‐Ping%Zhou%NULL%1,  Caspar W. N.%Looman%NULL%1,  Jan Hendrik%Richardus%NULL%1,  Wu‐Chun%Cao%NULL%1]"
$ws.Range("I7").Value = "_PMC"

# Row 8
$ws.Range("E8").Value = "[Mark%Loeb%NULL%1,  Allison%McGeer%NULL%2,  Bonnie%Henry%NULL%2,  Marianna%Ofner%NULL%2,  David%Rose%NULL%1,  Tammy%Hlywka%NULL%1,  Joanne%Levie%NULL%1,  Jane%McQueen%NULL%1,  Stephanie%Smith%NULL%1,  Lorraine%Moss%NULL%1,  Andrew%Smith%NULL%1,  Karen%Green%NULL%0,  Stephen D.%Walter%NULL%1]"
$ws.Range("I8").Value = "_PMC"

# Row 9
$ws.Range("E9").Value = "[Janet%Raboud%NULL%1,  Altynay%Shigayeva%NULL%1,  Allison%McGeer%NULL%0,  Erika%Bontovics%NULL%1,  Martin%Chapman%NULL%1,  Denise%Gravel%NULL%1,  Bonnie%Henry%NULL%0,  Stephen%Lapinsky%NULL%1,  Mark%Loeb%NULL%1,  L. Clifford%McDonald%NULL%0,  Marianna%Ofner%NULL%0,  Shirley%Paton%NULL%1,  Donna%Reynolds%NULL%1,  Damon%Scales%NULL%1,  Sandy%Shen%NULL%1,  Andrew%Simor%NULL%1,  Thomas%Stewart%NULL%1,  Mary%Vearncombe%NULL%1,  Dick%Zoutman%NULL%1,  Karen%Green%NULL%0,  Joel Mark%Montgomery%NULL%4,  Joel Mark%Montgomery%NULL%0]"
$ws.Range("I9").Value = "_PMC"

# Row 10 - only adds I10 (empty string)
$ws.Range("I10").Value = ""

# Row 11
$ws.Range("E11").Value = "[Hyungoo%Shin%NULL%1,  Jaehoon%Oh%NULL%1,  Tae Ho%Lim%NULL%1,  Hyunggoo%Kang%NULL%1,  Yeongtak%Song%NULL%1,  Sanghyun%Lee%NULL%1,  Abdelouahab%Bellou.%NULL%2,  Abdelouahab%Bellou.%NULL%0]"
$ws.Range("I11").Value = "_PMC"

# Row 12 - only adds I12 (empty string)
$ws.Range("I12").Value = ""
